# Apply the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.771.10'
$ws.Range("E2").Value = '  -0.36%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.650.82'
$ws.Range("E3").Value = '  +0.27%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.99'
$ws.Range("E5").Value = '  -0.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '327.83'
$ws.Range("E6").Value = '  +0.43%  '

# Row 7
$ws.Range("E7").Value = '  -0.91%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.551'
$ws.Range("E9").Value = '  -0.83%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.91'
$ws.Range("E10").Value = '  -2.35%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.06'
$ws.Range("E11").Value = '  -0.27%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0817'
$ws.Range("E12").Value = '  -0.61%  '

# Row 13
$ws.Range("E13").Value = '  +2.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.60'
$ws.Range("E14").Value = '  +3.06%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.066.41'
$ws.Range("E15").Value = '  +0.15%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.654.26'
$ws.Range("E16").Value = '  +0.40%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.861'
$ws.Range("E17").Value = '  -1.14%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.735.33'
$ws.Range("E18").Value = '  -0.34%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  +1.32%  '

# Row 20
$ws.Range("B20").Value = 'ImmutableX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.92'
$ws.Range("E20").Value = '  +0.16%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.70'
$ws.Range("E21").Value = '  -0.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0952'
$ws.Range("E22").Value = '  -0.26%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.12'
$ws.Range("E23").Value = '  -2.37%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.32'
$ws.Range("E24").Value = '  -4.01%  '

# Row 25
$ws.Range("E25").Value = '  -0.42%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.21'
$ws.Range("E26").Value = '  -2.31%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.21'
$ws.Range("E28").Value = '  +1.91%  '

# Row 29
$ws.Range("E29").Value = '  -0.65%  '

# Row 30
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.138'
$ws.Range("E30").Value = '  -2.29%  '

# Row 31
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.99'
$ws.Range("E31").Value = '  -3.85%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.58'
$ws.Range("E32").Value = '  -1.13%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.49'
$ws.Range("E33").Value = '  +0.68%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0821'
$ws.Range("E34").Value = '  +0.55%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.22'
$ws.Range("E35").Value = '  -1.56%  '

# Row 36
$ws.Range("E36").Value = '  -0.18%  '

# Row 37
$ws.Range("E37").Value = '  -1.17%  '

# Row 38
$ws.Range("E38").Value = '  -1.14%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.13'
$ws.Range("E39").Value = '  +0.87%  '

# Row 40
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '128.95'
$ws.Range("E40").Value = '  +3.82%  '

# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.69'
$ws.Range("E41").Value = '  +7.72%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0345'
$ws.Range("E42").Value = '  +8.75%  '

# Row 43
$ws.Range("E43").Value = '  +2.19%  '

# Row 44
$ws.Range("E44").Value = '  -0.39%  '

# Row 45
$ws.Range("E45").Value = '  +0.05%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.068.01'
$ws.Range("E46").Value = '  -0.92%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.14'
$ws.Range("E47").Value = '  +7.37%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.19'
$ws.Range("E48").Value = '  -2.99%  '

# Row 49
$ws.Range("E49").Value = '  -2.13%  '

# Row 50
$ws.Range("E50").Value = '  -1.92%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.23'
$ws.Range("E51").Value = '  -0.71%  '
